$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 12
$ws.Cells.Item(12, 8).Value = 12688.375
$ws.Cells.Item(12, 9).Value = 16866.166
$ws.Cells.Item(12, 10).Value = 155
$ws.Cells.Item(12, 11).Value = 16866.166
$ws.Cells.Item(12, 12).Value = 155
$ws.Cells.Item(12, 13).Value = -16696.166
$ws.Cells.Item(12, 14).Value = -495
# Row 17
$ws.Cells.Item(17, 8).Value = 204389.84
$ws.Cells.Item(17, 10).Value = 214209.12
$ws.Cells.Item(17, 12).Value = 642627.36
$ws.Cells.Item(17, 14).Value = -642963.36
# Row 29
$ws.Cells.Item(29, 8).Value = 5219.8
$ws.Cells.Item(29, 9).Value = 0
$ws.Cells.Item(29, 11).Value = 0
$ws.Cells.Item(29, 13).ClearContents()
# Row 38
$ws.Cells.Item(38, 8).Value = 2422.4546
$ws.Cells.Item(38, 9).Value = 893.5
$ws.Cells.Item(38, 11).Value = 2680.5
$ws.Cells.Item(38, 13).Value = -2308.5
# Row 40
$ws.Cells.Item(40, 8).Value = 4254.8887
$ws.Cells.Item(40, 9).Value = 3825
$ws.Cells.Item(40, 11).Value = 3825
$ws.Cells.Item(40, 13).Value = -3650
# Row 41
$ws.Cells.Item(41, 8).Value = 2595.6
$ws.Cells.Item(41, 9).Value = 2869.75
$ws.Cells.Item(41, 10).Value = 1499
$ws.Cells.Item(41, 11).Value = 2869.75
$ws.Cells.Item(41, 12).Value = 1499
$ws.Cells.Item(41, 13).Value = -2429.75
$ws.Cells.Item(41, 14).Value = -2379
# Row 45
$ws.Cells.Item(45, 8).Value = 350
$ws.Cells.Item(45, 9).Value = 350
$ws.Cells.Item(45, 11).Value = 1050
$ws.Cells.Item(45, 13).Value = -858
# Row 70
$ws.Cells.Item(70, 8).Value = 4744.7393
$ws.Cells.Item(70, 9).Value = 2897.2307
$ws.Cells.Item(70, 10).Value = 7146.5
$ws.Cells.Item(70, 11).Value = 8691.6921
$ws.Cells.Item(70, 12).Value = 21439.5
$ws.Cells.Item(70, 13).Value = -8421.6921
$ws.Cells.Item(70, 14).Value = -21979.5
# Row 73
$ws.Cells.Item(73, 8).Value = 4744.7393
$ws.Cells.Item(73, 9).Value = 2897.2307
$ws.Cells.Item(73, 10).Value = 7146.5
$ws.Cells.Item(73, 11).Value = 8691.6921
$ws.Cells.Item(73, 12).Value = 21439.5
$ws.Cells.Item(73, 13).Value = -7755.6921
$ws.Cells.Item(73, 14).Value = -23311.5
# Row 87
$ws.Cells.Item(87, 8).Value = 36850.5
$ws.Cells.Item(87, 10).Value = 39134
$ws.Cells.Item(87, 12).Value = 39134
$ws.Cells.Item(87, 14).Value = -41630
# Row 88
$ws.Cells.Item(88, 8).Value = 24617774
$ws.Cells.Item(88, 9).Value = 66678268
$ws.Cells.Item(88, 10).Value = 3587527.2
$ws.Cells.Item(88, 11).Value = 66678268
$ws.Cells.Item(88, 12).Value = 3587527.2
$ws.Cells.Item(88, 13).Value = -66677862
$ws.Cells.Item(88, 14).Value = -3588339.2
# Row 90
$ws.Cells.Item(90, 8).Value = 36850.5
$ws.Cells.Item(90, 10).Value = 39134
$ws.Cells.Item(90, 12).Value = 117402
$ws.Cells.Item(90, 14).Value = -129882
# Row 91
$ws.Cells.Item(91, 8).Value = 24617774
$ws.Cells.Item(91, 9).Value = 66678268
$ws.Cells.Item(91, 10).Value = 3587527.2
$ws.Cells.Item(91, 11).Value = 66678268
$ws.Cells.Item(91, 12).Value = 3587527.2
$ws.Cells.Item(91, 13).Value = -66676864
$ws.Cells.Item(91, 14).Value = -3590335.2
# Row 92
$ws.Cells.Item(92, 8).Value = 814.26086
$ws.Cells.Item(92, 9).Value = 825.8421
$ws.Cells.Item(92, 11).Value = 825.8421
$ws.Cells.Item(92, 13).Value = 422.1579

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Cells.Item(45, 8).Value = 5465.44
$ws.Cells.Item(45, 9).Value = 5573.136
$ws.Cells.Item(45, 11).Value = 5573.136
$ws.Cells.Item(45, 13).Value = -5196.136
# Row 61
$ws.Cells.Item(61, 8).Value = 58826308
$ws.Cells.Item(61, 9).Value = 66669550
$ws.Cells.Item(61, 11).Value = 66669550
$ws.Cells.Item(61, 13).Value = -66669338
# Row 74
$ws.Cells.Item(74, 8).Value = 52635556
$ws.Cells.Item(74, 9).Value = 58827790
$ws.Cells.Item(74, 11).Value = 58827790
$ws.Cells.Item(74, 13).Value = -58826916
# Row 77
$ws.Cells.Item(77, 8).Value = 52635556
$ws.Cells.Item(77, 9).Value = 58827790
$ws.Cells.Item(77, 11).Value = 294138950
$ws.Cells.Item(77, 13).Value = -294134582
# Row 136
$ws.Cells.Item(136, 8).Value = 58826308
$ws.Cells.Item(136, 9).Value = 66669550
$ws.Cells.Item(136, 11).Value = 200008650
$ws.Cells.Item(136, 13).Value = -200006100

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Cells.Item(94, 8).Value = 5972.615
$ws.Cells.Item(94, 9).Value = 9025.429
$ws.Cells.Item(94, 11).Value = 9025.429
$ws.Cells.Item(94, 13).Value = -8574.429
# Row 99
$ws.Cells.Item(99, 8).Value = 2187.375
$ws.Cells.Item(99, 9).Value = 1899.9231
$ws.Cells.Item(99, 11).Value = 1899.9231
$ws.Cells.Item(99, 13).Value = -401.9231

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Cells.Item(7, 8).Value = 14746.25
$ws.Cells.Item(7, 9).Value = 19625
$ws.Cells.Item(7, 10).Value = 110
$ws.Cells.Item(7, 11).Value = 19625
$ws.Cells.Item(7, 12).Value = 110
$ws.Cells.Item(7, 13).Value = -19512
$ws.Cells.Item(7, 14).Value = -336
# Row 20
$ws.Cells.Item(20, 8).Value = 79500
$ws.Cells.Item(20, 10).Value = 79500
$ws.Cells.Item(20, 12).Value = 79500
$ws.Cells.Item(20, 14).Value = -79972
# Row 30
$ws.Cells.Item(30, 8).Value = 79500
$ws.Cells.Item(30, 10).Value = 79500
$ws.Cells.Item(30, 12).Value = 79500
$ws.Cells.Item(30, 14).Value = -79682
# Row 128
$ws.Cells.Item(128, 8).Value = 79500
$ws.Cells.Item(128, 10).Value = 79500
$ws.Cells.Item(128, 12).Value = 79500
$ws.Cells.Item(128, 14).Value = -89460
# Row 129
$ws.Cells.Item(129, 8).Value = 91600
$ws.Cells.Item(129, 9).Value = 85000
$ws.Cells.Item(129, 10).Value = 96000
$ws.Cells.Item(129, 11).Value = 85000
$ws.Cells.Item(129, 12).Value = 96000
$ws.Cells.Item(129, 13).Value = -80000
$ws.Cells.Item(129, 14).Value = -106000
# Row 130
$ws.Cells.Item(130, 8).Value = 89999
$ws.Cells.Item(130, 10).Value = 89999
$ws.Cells.Item(130, 12).Value = 89999
$ws.Cells.Item(130, 14).Value = -100039
# Row 131
$ws.Cells.Item(131, 8).Value = 92997
$ws.Cells.Item(131, 10).Value = 92997
$ws.Cells.Item(131, 12).Value = 92997
$ws.Cells.Item(131, 14).Value = -103077

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 12
$ws.Cells.Item(12, 8).Value = 220.76471
$ws.Cells.Item(12, 10).Value = 311.625
$ws.Cells.Item(12, 12).Value = 934.875
$ws.Cells.Item(12, 14).Value = -1280.875
# Row 107
$ws.Cells.Item(107, 8).Value = 1417.1818
$ws.Cells.Item(107, 10).Value = 1766.875
$ws.Cells.Item(107, 12).Value = 5300.625
$ws.Cells.Item(107, 14).Value = -9140.625

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Cells.Item(2, 8).Value = 281.77777
$ws.Cells.Item(2, 9).Value = 260
$ws.Cells.Item(2, 10).Value = 288
$ws.Cells.Item(2, 11).Value = 260
$ws.Cells.Item(2, 12).Value = 288
$ws.Cells.Item(2, 13).Value = -147
$ws.Cells.Item(2, 14).Value = -514
# Row 132
$ws.Cells.Item(132, 8).Value = 6256382.5
$ws.Cells.Item(132, 9).Value = 7359753
$ws.Cells.Item(132, 11).Value = 22079259
$ws.Cells.Item(132, 13).Value = -22076729

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Cells.Item(7, 8).Value = 4591.125
$ws.Cells.Item(7, 9).Value = 4591.125
$ws.Cells.Item(7, 11).Value = 4591.125
$ws.Cells.Item(7, 13).Value = -4479.125
# Row 68
$ws.Cells.Item(68, 8).Value = 2070.3333
$ws.Cells.Item(68, 9).Value = 1922.875
$ws.Cells.Item(68, 11).Value = 1922.875
$ws.Cells.Item(68, 13).Value = -1173.875
# Row 71
$ws.Cells.Item(71, 8).Value = 2070.3333
$ws.Cells.Item(71, 9).Value = 1922.875
$ws.Cells.Item(71, 11).Value = 9614.375
$ws.Cells.Item(71, 13).Value = -5870.375
# Row 100
$ws.Cells.Item(100, 8).Value = 7678036.5
$ws.Cells.Item(100, 9).Value = 8679337
$ws.Cells.Item(100, 10).Value = 1399.6666
$ws.Cells.Item(100, 11).Value = 8679337
$ws.Cells.Item(100, 12).Value = 1399.6666
$ws.Cells.Item(100, 13).Value = -8678796
$ws.Cells.Item(100, 14).Value = -2481.6666
# Row 126
$ws.Cells.Item(126, 8).Value = 4591.125
$ws.Cells.Item(126, 9).Value = 4591.125
$ws.Cells.Item(126, 11).Value = 13773.375
$ws.Cells.Item(126, 13).Value = -11303.375
